$d = $word.ActiveDocument

$pairs = @(
    @("44×76=", "67×65="),
    @("65×18=", "76×60="),
    @("80×56=", "84×70="),
    @("70×63=", "18×44="),
    @("77×48=", "71×43="),
    @("67×27=", "80×91="),
    @("82×28=", "60×97="),
    @("28×39=", "82×82="),
    @("55×68=", "79×82="),
    @("91×84=", "87×43="),
    @("44×47=", "23×32="),
    @("68×54=", "83×87="),
    @("83×53=", "69×58="),
    @("98×46=", "83×76="),
    @("79×21=", "11×20="),
    @("29×94=", "21×14="),
    @("22×49=", "11×52="),
    @("39×66=", "31×22="),
    @("65×22=", "51×46="),
    @("82×40=", "57×56="),
    @("47×47=", "64×29="),
    @("27×19=", "41×52="),
    @("36×97=", "97×15="),
    @("53×21=", "47×98="),
    @("26×79=", "93×11=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
